$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# Update the cell content: "Good Morning" -> "GIT UPDATE"
$ws.Range("E8").Value = "GIT UPDATE"

# Reflect the active selection recorded in the sheet view (activeCell="E8" sqref="E8")
$ws.Activate()
$ws.Range("E8").Select()
